$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Reword existing range-of-function-values descriptions (B167:B170)
# and add the 10 new rows (172-181) for the new building-block pairs (prob30, p&s part).

$ws.Range("A166").Value = 'x0016'
$ws.Range("B166").Value = '함수의 조건에 맞는 가능한 치역을 모두 구합니다.'

$ws.Range("A167").Value = 'x0017'
$ws.Range("B167").Value = '치역이 $\{1, 2, 3\}$일 때, 조건에 맞는 함수의 개수를 구합니다.'
$ws.Range("C167").Value = '32111_x28'

$ws.Range("A168").Value = 'x0018'
$ws.Range("B168").Value = '치역이 $\{1, 2, 4\}$일 때, 조건에 맞는 함수의 개수를 구합니다.'
$ws.Range("C168").Value = '32111_x28'

$ws.Range("A169").Value = 'x0019'
$ws.Range("B169").Value = '치역이 $\{1, 3, 4\}$일 때, 조건에 맞는 함수의 개수를 구합니다.'
$ws.Range("C169").Value = '32111_x28'

$ws.Range("A170").Value = 'x0020'
$ws.Range("B170").Value = '치역이 $\{2, 3, 4\}$일 때, 조건에 맞는 함수의 개수를 구합니다.'
$ws.Range("C170").Value = '32111_x28'

$ws.Range("A171").Value = 'x0021'
$ws.Range("B171").Value = '각각의 개수를 모두 더해서 조건을 만족시키는 전체 개수를 구합니다. '

$ws.Range("A172").Value = 'x0022'
$ws.Range("B172").Value = '연속확률변수의 확률밀도함수의 함숫값이 $0$ 이상이 되도록하는 조건을 구합니다'

$ws.Range("A173").Value = 'x0023'
$ws.Range("B173").Value = '연속확률변수의 확률밀도함수로 둘러싸인 넓이가 $1$이 되도록하는 조건을 합니다.'

$ws.Range("A174").Value = 'x0024'
$ws.Range("B174").Value = '확률밀도함수가 되도록 방정식을 세워 확률밀도함수의 미정계수를 정합니다.'

$ws.Range("A175").Value = 'x0025'
$ws.Range("B175").Value = '연속확률변수의 주어진 범위에서의 확률에 해당하는 확률밀도함수와 둘러싸인 도형의 넓이를 구합니다.'

$ws.Range("A176").Value = 'x0026'
$ws.Range("B176").Value = '$a_{5}+b_{5} \geq 7$인 사건을 합 $a_{5}+b_{5}$에 대해 분류합니다.'
$ws.Range("C176").Value = '32111_x30'

$ws.Range("A177").Value = 'x0027'
$ws.Range("B177").Value = '합 $a_{5}+b_{5}$이 $7$, $8$, $9$, $10$인 독립시행의 확률을 각각 구합니다.'
$ws.Range("C177").Value = '32111_x30'

$ws.Range("A178").Value = 'x0028'
$ws.Range("B178").Value = '각 경우의 확률을 모두 더해 $a_{5}+b_{5} \geq 7$인 사건의 확률을 구합니다.'
$ws.Range("C178").Value = '32111_x30'

$ws.Range("A179").Value = 'x0029'
$ws.Range("B179").Value = '합 $a_{5}+b_{5}$이 $7$, $8$, $9$, $10$인 각 경우에 대해 $a_{k}=b_{k}$ 인 자연수 $1 \leq k \leq 5$가 존재할 확률을 구합니다.'
$ws.Range("C179").Value = '32111_x30'

$ws.Range("A180").Value = 'x0030'
$ws.Range("B180").Value = '각 경우의 확률을 모두 더해 $a_{5}+b_{5} \geq 7$이고, $a_{k}=b_{k}$ 인 자연수 $1 \leq k \leq 5$가 존재할 확률울 구합니다.'
$ws.Range("C180").Value = '32111_x30'

$ws.Range("A181").Value = 'x0031'
$ws.Range("B181").Value = '구해진 두 확률의 비를 통해서 문제에서 요구하는 조건부 확률을 구합니다. '

# Old footer marker rows (192, 202) shift up to (191, 201) to follow the newly inserted rows
$ws.Range("A191").Value = $ws.Range("A192").Value()
$ws.Range("A192").ClearContents()
$ws.Range("A201").Value = $ws.Range("A202").Value()
$ws.Range("A202").ClearContents()

# Restore the view state (active cell selection) to reflect the new layout
$ws.Activate() | Out-Null
$ws.Range("B182").Select() | Out-Null
